$wb = $excel.ActiveWorkbook

$changes = @(
    @{ Sheet = "ALC"; Addr = "H33"; Value = 540.6875 },
    @{ Sheet = "ALC"; Addr = "I33"; Value = 329.5 },
    @{ Sheet = "ALC"; Addr = "J33"; Value = 1174.25 },
    @{ Sheet = "ALC"; Addr = "K33"; Value = 329.5 },
    @{ Sheet = "ALC"; Addr = "L33"; Value = 1174.25 },
    @{ Sheet = "ALC"; Addr = "M33"; Value = -100.5 },
    @{ Sheet = "ALC"; Addr = "N33"; Value = -1632.25 },
    @{ Sheet = "ALC"; Addr = "H51"; Value = 11774.625 },
    @{ Sheet = "ALC"; Addr = "I51"; Value = 9799 },
    @{ Sheet = "ALC"; Addr = "K51"; Value = 9799 },
    @{ Sheet = "ALC"; Addr = "M51"; Value = -9315 },
    @{ Sheet = "ALC"; Addr = "H137"; Value = 4131 },
    @{ Sheet = "ALC"; Addr = "I137"; Value = 2566 },
    @{ Sheet = "ALC"; Addr = "J137"; Value = 4717.875 },
    @{ Sheet = "ALC"; Addr = "K137"; Value = 7698 },
    @{ Sheet = "ALC"; Addr = "L137"; Value = 14153.625 },
    @{ Sheet = "ALC"; Addr = "M137"; Value = -5148 },
    @{ Sheet = "ALC"; Addr = "N137"; Value = -19253.625 },
    @{ Sheet = "ARM"; Addr = "H61"; Value = 4101.524 },
    @{ Sheet = "ARM"; Addr = "I61"; Value = 2849.5293 },
    @{ Sheet = "ARM"; Addr = "J61"; Value = 4564.2173 },
    @{ Sheet = "ARM"; Addr = "K61"; Value = 2849.5293 },
    @{ Sheet = "ARM"; Addr = "L61"; Value = 4564.2173 },
    @{ Sheet = "ARM"; Addr = "M61"; Value = -2637.5293 },
    @{ Sheet = "ARM"; Addr = "N61"; Value = -4988.2173 },
    @{ Sheet = "ARM"; Addr = "H74"; Value = 3181.1052 },
    @{ Sheet = "ARM"; Addr = "I74"; Value = 1690.25 },
    @{ Sheet = "ARM"; Addr = "J74"; Value = 4265.364 },
    @{ Sheet = "ARM"; Addr = "K74"; Value = 1690.25 },
    @{ Sheet = "ARM"; Addr = "L74"; Value = 4265.364 },
    @{ Sheet = "ARM"; Addr = "M74"; Value = -816.25 },
    @{ Sheet = "ARM"; Addr = "N74"; Value = -6013.364 },
    @{ Sheet = "ARM"; Addr = "H77"; Value = 3181.1052 },
    @{ Sheet = "ARM"; Addr = "I77"; Value = 1690.25 },
    @{ Sheet = "ARM"; Addr = "J77"; Value = 4265.364 },
    @{ Sheet = "ARM"; Addr = "K77"; Value = 8451.25 },
    @{ Sheet = "ARM"; Addr = "L77"; Value = 21326.82 },
    @{ Sheet = "ARM"; Addr = "M77"; Value = -4083.25 },
    @{ Sheet = "ARM"; Addr = "N77"; Value = -30062.82 },
    @{ Sheet = "ARM"; Addr = "H102"; Value = 2315.353 },
    @{ Sheet = "ARM"; Addr = "I102"; Value = 2279.4546 },
    @{ Sheet = "ARM"; Addr = "K102"; Value = 2279.4546 },
    @{ Sheet = "ARM"; Addr = "M102"; Value = -657.4546 },
    @{ Sheet = "ARM"; Addr = "H132"; Value = 1905586.6 },
    @{ Sheet = "ARM"; Addr = "I132"; Value = 2650492.8 },
    @{ Sheet = "ARM"; Addr = "J132"; Value = 167472.33 },
    @{ Sheet = "ARM"; Addr = "K132"; Value = 7951478.399999999 },
    @{ Sheet = "ARM"; Addr = "L132"; Value = 502416.99 },
    @{ Sheet = "ARM"; Addr = "M132"; Value = -7948948.399999999 },
    @{ Sheet = "ARM"; Addr = "N132"; Value = -507476.99 },
    @{ Sheet = "ARM"; Addr = "H136"; Value = 4101.524 },
    @{ Sheet = "ARM"; Addr = "I136"; Value = 2849.5293 },
    @{ Sheet = "ARM"; Addr = "J136"; Value = 4564.2173 },
    @{ Sheet = "ARM"; Addr = "K136"; Value = 8548.5879 },
    @{ Sheet = "ARM"; Addr = "L136"; Value = 13692.6519 },
    @{ Sheet = "ARM"; Addr = "M136"; Value = -5998.5879 },
    @{ Sheet = "ARM"; Addr = "N136"; Value = -18792.6519 },
    @{ Sheet = "BSM"; Addr = "H80"; Value = 18533240 },
    @{ Sheet = "BSM"; Addr = "I80"; Value = 1425.5 },
    @{ Sheet = "BSM"; Addr = "J80"; Value = 33358694 },
    @{ Sheet = "BSM"; Addr = "K80"; Value = 1425.5 },
    @{ Sheet = "BSM"; Addr = "L80"; Value = 33358694 },
    @{ Sheet = "BSM"; Addr = "M80"; Value = -427.5 },
    @{ Sheet = "BSM"; Addr = "N80"; Value = -33360690 },
    @{ Sheet = "BSM"; Addr = "H83"; Value = 18533240 },
    @{ Sheet = "BSM"; Addr = "I83"; Value = 1425.5 },
    @{ Sheet = "BSM"; Addr = "J83"; Value = 33358694 },
    @{ Sheet = "BSM"; Addr = "K83"; Value = 7127.5 },
    @{ Sheet = "BSM"; Addr = "L83"; Value = 166793470 },
    @{ Sheet = "BSM"; Addr = "M83"; Value = -2135.5 },
    @{ Sheet = "BSM"; Addr = "N83"; Value = -166803454 },
    @{ Sheet = "BSM"; Addr = "H134"; Value = 5977.278 },
    @{ Sheet = "BSM"; Addr = "I134"; Value = 4793.9287 },
    @{ Sheet = "BSM"; Addr = "J134"; Value = 10119 },
    @{ Sheet = "BSM"; Addr = "K134"; Value = 14381.7861 },
    @{ Sheet = "BSM"; Addr = "L134"; Value = 30357 },
    @{ Sheet = "BSM"; Addr = "M134"; Value = -11846.7861 },
    @{ Sheet = "BSM"; Addr = "N134"; Value = -35427 },
    @{ Sheet = "CRP"; Addr = "H7"; Value = 108.2 },
    @{ Sheet = "CRP"; Addr = "I7"; Value = 129 },
    @{ Sheet = "CRP"; Addr = "K7"; Value = 129 },
    @{ Sheet = "CRP"; Addr = "M7"; Value = -16 },
    @{ Sheet = "CRP"; Addr = "H22"; Value = 851552.5 },
    @{ Sheet = "CRP"; Addr = "I22"; Value = 1623733.9 },
    @{ Sheet = "CRP"; Addr = "K22"; Value = 1623733.9 },
    @{ Sheet = "CRP"; Addr = "M22"; Value = -1623383.9 },
    @{ Sheet = "CRP"; Addr = "H31"; Value = 11117934 },
    @{ Sheet = "CRP"; Addr = "I31"; Value = 45473068 },
    @{ Sheet = "CRP"; Addr = "J31"; Value = 3037.2646 },
    @{ Sheet = "CRP"; Addr = "K31"; Value = 45473068 },
    @{ Sheet = "CRP"; Addr = "L31"; Value = 3037.2646 },
    @{ Sheet = "CRP"; Addr = "M31"; Value = -45472773 },
    @{ Sheet = "CRP"; Addr = "N31"; Value = -3627.2646 },
    @{ Sheet = "CRP"; Addr = "H34"; Value = 11117934 },
    @{ Sheet = "CRP"; Addr = "I34"; Value = 45473068 },
    @{ Sheet = "CRP"; Addr = "J34"; Value = 3037.2646 },
    @{ Sheet = "CRP"; Addr = "K34"; Value = 45473068 },
    @{ Sheet = "CRP"; Addr = "L34"; Value = 3037.2646 },
    @{ Sheet = "CRP"; Addr = "M34"; Value = -45472866 },
    @{ Sheet = "CRP"; Addr = "N34"; Value = -3441.2646 },
    @{ Sheet = "CRP"; Addr = "H58"; Value = 52637400 },
    @{ Sheet = "CRP"; Addr = "I58"; Value = 100002890 },
    @{ Sheet = "CRP"; Addr = "J58"; Value = 9082 },
    @{ Sheet = "CRP"; Addr = "K58"; Value = 100002890 },
    @{ Sheet = "CRP"; Addr = "L58"; Value = 9082 },
    @{ Sheet = "CRP"; Addr = "M58"; Value = -100002687 },
    @{ Sheet = "CRP"; Addr = "N58"; Value = -9488 },
    @{ Sheet = "CRP"; Addr = "H99"; Value = 5559202 },
    @{ Sheet = "CRP"; Addr = "I99"; Value = 9263110 },
    @{ Sheet = "CRP"; Addr = "J99"; Value = 3338.875 },
    @{ Sheet = "CRP"; Addr = "K99"; Value = 9263110 },
    @{ Sheet = "CRP"; Addr = "L99"; Value = 3338.875 },
    @{ Sheet = "CRP"; Addr = "M99"; Value = -9261612 },
    @{ Sheet = "CRP"; Addr = "N99"; Value = -6334.875 },
    @{ Sheet = "CRP"; Addr = "H105"; Value = 62501468 },
    @{ Sheet = "CRP"; Addr = "I105"; Value = 71429320 },
    @{ Sheet = "CRP"; Addr = "J105"; Value = 6500 },
    @{ Sheet = "CRP"; Addr = "K105"; Value = 71429320 },
    @{ Sheet = "CRP"; Addr = "L105"; Value = 6500 },
    @{ Sheet = "CRP"; Addr = "M105"; Value = -71427573 },
    @{ Sheet = "CRP"; Addr = "N105"; Value = -9994 },
    @{ Sheet = "CRP"; Addr = "H107"; Value = 618.8570999999999 },
    @{ Sheet = "CRP"; Addr = "I107"; Value = 276.9 },
    @{ Sheet = "CRP"; Addr = "J107"; Value = 1473.75 },
    @{ Sheet = "CRP"; Addr = "K107"; Value = 276.9 },
    @{ Sheet = "CRP"; Addr = "L107"; Value = 1473.75 },
    @{ Sheet = "CRP"; Addr = "M107"; Value = 1643.1 },
    @{ Sheet = "CRP"; Addr = "N107"; Value = -5313.75 },
    @{ Sheet = "CRP"; Addr = "H126"; Value = 5559202 },
    @{ Sheet = "CRP"; Addr = "I126"; Value = 9263110 },
    @{ Sheet = "CRP"; Addr = "J126"; Value = 3338.875 },
    @{ Sheet = "CRP"; Addr = "K126"; Value = 27789330 },
    @{ Sheet = "CRP"; Addr = "L126"; Value = 10016.625 },
    @{ Sheet = "CRP"; Addr = "M126"; Value = -27786860 },
    @{ Sheet = "CRP"; Addr = "N126"; Value = -14956.625 },
    @{ Sheet = "CRP"; Addr = "H132"; Value = 4060.8125 },
    @{ Sheet = "CRP"; Addr = "I132"; Value = 2635 },
    @{ Sheet = "CRP"; Addr = "J132"; Value = 14041.5 },
    @{ Sheet = "CRP"; Addr = "K132"; Value = 7905 },
    @{ Sheet = "CRP"; Addr = "L132"; Value = 42124.5 },
    @{ Sheet = "CRP"; Addr = "M132"; Value = -5375 },
    @{ Sheet = "CRP"; Addr = "N132"; Value = -47184.5 },
    @{ Sheet = "CRP"; Addr = "H134"; Value = 103900980 },
    @{ Sheet = "CRP"; Addr = "I134"; Value = 126986216 },
    @{ Sheet = "CRP"; Addr = "J134"; Value = 17375 },
    @{ Sheet = "CRP"; Addr = "K134"; Value = 380958648 },
    @{ Sheet = "CRP"; Addr = "L134"; Value = 52125 },
    @{ Sheet = "CRP"; Addr = "M134"; Value = -380956113 },
    @{ Sheet = "CRP"; Addr = "N134"; Value = -57195 },
    @{ Sheet = "CRP"; Addr = "H136"; Value = 52637400 },
    @{ Sheet = "CRP"; Addr = "I136"; Value = 100002890 },
    @{ Sheet = "CRP"; Addr = "J136"; Value = 9082 },
    @{ Sheet = "CRP"; Addr = "K136"; Value = 300008670 },
    @{ Sheet = "CRP"; Addr = "L136"; Value = 27246 },
    @{ Sheet = "CRP"; Addr = "M136"; Value = -300006120 },
    @{ Sheet = "CRP"; Addr = "N136"; Value = -32346 },
    @{ Sheet = "CUL"; Addr = "H14"; Value = 775.5714 },
    @{ Sheet = "CUL"; Addr = "I14"; Value = 775.5714 },
    @{ Sheet = "CUL"; Addr = "K14"; Value = 2326.7142 },
    @{ Sheet = "CUL"; Addr = "M14"; Value = -2153.7142 },
    @{ Sheet = "CUL"; Addr = "H33"; Value = 697.375 },
    @{ Sheet = "CUL"; Addr = "J33"; Value = 976.2 },
    @{ Sheet = "CUL"; Addr = "L33"; Value = 5857.200000000001 },
    @{ Sheet = "CUL"; Addr = "N33"; Value = -6423.200000000001 },
    @{ Sheet = "CUL"; Addr = "H80"; Value = 1749.375 },
    @{ Sheet = "CUL"; Addr = "J80"; Value = 1665.8334 },
    @{ Sheet = "CUL"; Addr = "L80"; Value = 4997.5002 },
    @{ Sheet = "CUL"; Addr = "N80"; Value = -6869.5002 },
    @{ Sheet = "CUL"; Addr = "H83"; Value = 1749.375 },
    @{ Sheet = "CUL"; Addr = "J83"; Value = 1665.8334 },
    @{ Sheet = "CUL"; Addr = "L83"; Value = 14992.5006 },
    @{ Sheet = "CUL"; Addr = "N83"; Value = -24352.5006 },
    @{ Sheet = "CUL"; Addr = "H86"; Value = 1755 },
    @{ Sheet = "CUL"; Addr = "J86"; Value = 380.83334 },
    @{ Sheet = "CUL"; Addr = "L86"; Value = 1142.50002 },
    @{ Sheet = "CUL"; Addr = "N86"; Value = -3514.50002 },
    @{ Sheet = "CUL"; Addr = "H89"; Value = 1755 },
    @{ Sheet = "CUL"; Addr = "J89"; Value = 380.83334 },
    @{ Sheet = "CUL"; Addr = "L89"; Value = 3427.50006 },
    @{ Sheet = "CUL"; Addr = "N89"; Value = -15283.50006 },
    @{ Sheet = "CUL"; Addr = "H98"; Value = 5799.2 },
    @{ Sheet = "CUL"; Addr = "I98"; Value = 1499.5 },
    @{ Sheet = "CUL"; Addr = "J98"; Value = 8665.666999999999 },
    @{ Sheet = "CUL"; Addr = "K98"; Value = 4498.5 },
    @{ Sheet = "CUL"; Addr = "L98"; Value = 25997.001 },
    @{ Sheet = "CUL"; Addr = "M98"; Value = -3000.5 },
    @{ Sheet = "CUL"; Addr = "N98"; Value = -28993.001 },
    @{ Sheet = "CUL"; Addr = "H131"; Value = 14289262 },
    @{ Sheet = "CUL"; Addr = "I131"; Value = 22226504 },
    @{ Sheet = "CUL"; Addr = "J131"; Value = 2228.8 },
    @{ Sheet = "CUL"; Addr = "K131"; Value = 66679512 },
    @{ Sheet = "CUL"; Addr = "L131"; Value = 6686.400000000001 },
    @{ Sheet = "CUL"; Addr = "M131"; Value = -66674472 },
    @{ Sheet = "CUL"; Addr = "N131"; Value = -16766.4 },
    @{ Sheet = "CUL"; Addr = "H141"; Value = 5095.8 },
    @{ Sheet = "CUL"; Addr = "I141"; Value = 5095.8 },
    @{ Sheet = "CUL"; Addr = "K141"; Value = 15287.4 },
    @{ Sheet = "CUL"; Addr = "M141"; Value = -10107.4 },
    @{ Sheet = "GSM"; Addr = "H80"; Value = 5333.222 },
    @{ Sheet = "GSM"; Addr = "I80"; Value = 2607.1428 },
    @{ Sheet = "GSM"; Addr = "J80"; Value = 14874.5 },
    @{ Sheet = "GSM"; Addr = "K80"; Value = 2607.1428 },
    @{ Sheet = "GSM"; Addr = "L80"; Value = 14874.5 },
    @{ Sheet = "GSM"; Addr = "M80"; Value = -1609.1428 },
    @{ Sheet = "GSM"; Addr = "N80"; Value = -16870.5 },
    @{ Sheet = "GSM"; Addr = "H83"; Value = 5333.222 },
    @{ Sheet = "GSM"; Addr = "I83"; Value = 2607.1428 },
    @{ Sheet = "GSM"; Addr = "J83"; Value = 14874.5 },
    @{ Sheet = "GSM"; Addr = "K83"; Value = 13035.714 },
    @{ Sheet = "GSM"; Addr = "L83"; Value = 74372.5 },
    @{ Sheet = "GSM"; Addr = "M83"; Value = -8043.714 },
    @{ Sheet = "GSM"; Addr = "N83"; Value = -84356.5 },
    @{ Sheet = "GSM"; Addr = "H107"; Value = 980.6667 },
    @{ Sheet = "GSM"; Addr = "I107"; Value = 185.75 },
    @{ Sheet = "GSM"; Addr = "J107"; Value = 1616.6 },
    @{ Sheet = "GSM"; Addr = "K107"; Value = 185.75 },
    @{ Sheet = "GSM"; Addr = "L107"; Value = 1616.6 },
    @{ Sheet = "GSM"; Addr = "M107"; Value = 1734.25 },
    @{ Sheet = "GSM"; Addr = "N107"; Value = -5456.6 },
    @{ Sheet = "GSM"; Addr = "H132"; Value = 3956.524 },
    @{ Sheet = "GSM"; Addr = "I132"; Value = 3812.6453 },
    @{ Sheet = "GSM"; Addr = "J132"; Value = 4362 },
    @{ Sheet = "GSM"; Addr = "K132"; Value = 11437.9359 },
    @{ Sheet = "GSM"; Addr = "L132"; Value = 13086 },
    @{ Sheet = "GSM"; Addr = "M132"; Value = -8907.9359 },
    @{ Sheet = "GSM"; Addr = "N132"; Value = -18146 },
    @{ Sheet = "LTW"; Addr = "H61"; Value = 7380.609 },
    @{ Sheet = "LTW"; Addr = "I61"; Value = 5478.294 },
    @{ Sheet = "LTW"; Addr = "K61"; Value = 5478.294 },
    @{ Sheet = "LTW"; Addr = "M61"; Value = -5276.294 },
    @{ Sheet = "LTW"; Addr = "H113"; Value = 7380.609 },
    @{ Sheet = "LTW"; Addr = "I113"; Value = 5478.294 },
    @{ Sheet = "LTW"; Addr = "K113"; Value = 5478.294 },
    @{ Sheet = "LTW"; Addr = "M113"; Value = -3308.294 },
    @{ Sheet = "LTW"; Addr = "H132"; Value = 7146519.5 },
    @{ Sheet = "LTW"; Addr = "I132"; Value = 10206128 },
    @{ Sheet = "LTW"; Addr = "K132"; Value = 30618384 },
    @{ Sheet = "LTW"; Addr = "M132"; Value = -30615854 },
    @{ Sheet = "WVR"; Addr = "H100"; Value = 1514.2858 },
    @{ Sheet = "WVR"; Addr = "I100"; Value = 1250.0667 },
    @{ Sheet = "WVR"; Addr = "J100"; Value = 2174.8333 },
    @{ Sheet = "WVR"; Addr = "K100"; Value = 2500.1334 },
    @{ Sheet = "WVR"; Addr = "L100"; Value = 4349.6666 },
    @{ Sheet = "WVR"; Addr = "M100"; Value = -1959.1334 },
    @{ Sheet = "WVR"; Addr = "N100"; Value = -5431.6666 },
    @{ Sheet = "WVR"; Addr = "H122"; Value = 5973.049 },
    @{ Sheet = "WVR"; Addr = "J122"; Value = 29999 },
    @{ Sheet = "WVR"; Addr = "L122"; Value = 89997 },
    @{ Sheet = "WVR"; Addr = "N122"; Value = -94897 },
    @{ Sheet = "WVR"; Addr = "H126"; Value = 4780.5 },
    @{ Sheet = "WVR"; Addr = "I126"; Value = 2050.1 },
    @{ Sheet = "WVR"; Addr = "J126"; Value = 8193.5 },
    @{ Sheet = "WVR"; Addr = "K126"; Value = 6150.299999999999 },
    @{ Sheet = "WVR"; Addr = "L126"; Value = 24580.5 },
    @{ Sheet = "WVR"; Addr = "M126"; Value = -3680.299999999999 },
    @{ Sheet = "WVR"; Addr = "N126"; Value = -29520.5 },
    @{ Sheet = "WVR"; Addr = "H132"; Value = 18400064 },
    @{ Sheet = "WVR"; Addr = "I132"; Value = 14496981 },
    @{ Sheet = "WVR"; Addr = "J132"; Value = 33361882 },
    @{ Sheet = "WVR"; Addr = "K132"; Value = 43490943 },
    @{ Sheet = "WVR"; Addr = "L132"; Value = 100085646 },
    @{ Sheet = "WVR"; Addr = "M132"; Value = -43488413 },
    @{ Sheet = "WVR"; Addr = "N132"; Value = -100090706 },
    @{ Sheet = "WVR"; Addr = "H136"; Value = 20845542 },
    @{ Sheet = "WVR"; Addr = "I136"; Value = 50022548 },
    @{ Sheet = "WVR"; Addr = "J136"; Value = 4822.9287 },
    @{ Sheet = "WVR"; Addr = "K136"; Value = 150067644 },
    @{ Sheet = "WVR"; Addr = "L136"; Value = 14468.7861 },
    @{ Sheet = "WVR"; Addr = "M136"; Value = -150065094 },
    @{ Sheet = "WVR"; Addr = "N136"; Value = -19568.7861 },
)

foreach ($ch in $changes) {
    $ws = $wb.Worksheets.Item($ch.Sheet)
    $ws.Range($ch.Addr).Value = $ch.Value
}

Write-Host "Applied $($changes.Count) cell updates."